$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.436241610738255
$ws.Range("C2").Value = 0.647727272727273
$ws.Range("D2").Value = 0.596774193548387
$ws.Range("E2").Value = 0.425287356321839
